$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.095.39'
$ws.Range("E2").Value = '  -0.25%  '
$ws.Range("D3").Value = '2.543.05'
$ws.Range("E3").Value = '  +4.52%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '568.98'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +0.60%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.63'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +4.83%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.582'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.86%  '
$ws.Range("D9").Value = '2.539.77'
$ws.Range("E9").Value = '  +4.52%  '
$ws.Range("E10").Value = '  -0.07%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.64'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.67%  '
$ws.Range("E12").Value = '  +0.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.355'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  +1.07%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.51'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  +4.55%  '
$ws.Range("D15").Value = '2.992.08'
$ws.Range("E15").Value = '  +4.34%  '
$ws.Range("D16").Value = '62.894.26'
$ws.Range("E16").Value = '  -0.29%  '
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("D18").Value = '2.535.64'
$ws.Range("E18").Value = '  +4.44%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '11.56'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +2.92%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '335.62'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -1.45%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.30'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  +1.54%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.77'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.12%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.25'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.02%  '
$ws.Range("E25").Value = '  -2.91%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.59'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +3.63%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.52'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  +12.28%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.997'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.35%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.38'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +2.46%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.21'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +8.02%  '
$ws.Range("D31").Value = '0.0₃0813'
$ws.Range("E31").Value = '  +1.49%  '
$ws.Range("E32").Value = '  +1.24%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '177.58'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +2.01%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.59'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +7.53%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '414.55'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +11.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.396'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -0.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '18.85'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  +1.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '4.41'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -1.73%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.76'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  +3.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.998'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -0.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '39.28'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -1.52%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '151.97'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  +3.17%  '
$ws.Range("E44").Value = '  +1.36%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '20.70'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  +0.53%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.607'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +2.32%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0967'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  +0.97%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0521'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  -0.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0236'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  +5.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.38'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +2.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.78'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +2.39%  '
